$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3666.5
$ws.Range("I64").Value = 3599.8
$ws.Range("K64").Value = 3599.8
$ws.Range("M64").Value = -3351.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3666.5
$ws.Range("I67").Value = 3599.8
$ws.Range("K67").Value = 3599.8
$ws.Range("M67").Value = -2741.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2515.1875
$ws.Range("I138").Value = 1974.3334
$ws.Range("J138").Value = 3210.5715
$ws.Range("K138").Value = 5923.0002
$ws.Range("L138").Value = 9631.7145
$ws.Range("M138").Value = -783.0002000000004
$ws.Range("N138").Value = -19911.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3969.9285
$ws.Range("I63").Value = 2506.8
$ws.Range("K63").Value = 2506.8
$ws.Range("M63").Value = -1820.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3969.9285
$ws.Range("I66").Value = 2506.8
$ws.Range("K66").Value = 12534
$ws.Range("M66").Value = -9102

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2360.8
$ws.Range("I74").Value = 1320.8125
$ws.Range("K74").Value = 1320.8125
$ws.Range("M74").Value = -446.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2360.8
$ws.Range("I77").Value = 1320.8125
$ws.Range("K77").Value = 6604.0625
$ws.Range("M77").Value = -2236.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 629
$ws.Range("I97").Value = 611.3570999999999
$ws.Range("J97").Value = 999.5
$ws.Range("K97").Value = 611.3570999999999
$ws.Range("L97").Value = 999.5
$ws.Range("M97").Value = -115.3570999999999
$ws.Range("N97").Value = -1991.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4177.727
$ws.Range("I102").Value = 4220.5
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 4220.5
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -2598.5
$ws.Range("N102").Value = -6994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15024.223
$ws.Range("J82").Value = 21855.334
$ws.Range("L82").Value = 21855.334
$ws.Range("N82").Value = -22621.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 15024.223
$ws.Range("J85").Value = 21855.334
$ws.Range("L85").Value = 21855.334
$ws.Range("N85").Value = -24507.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 436969.75
$ws.Range("I86").Value = 557594.9
$ws.Range("J86").Value = 2719.2
$ws.Range("K86").Value = 557594.9
$ws.Range("L86").Value = 2719.2
$ws.Range("M86").Value = -556471.9
$ws.Range("N86").Value = -4965.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 436969.75
$ws.Range("I89").Value = 557594.9
$ws.Range("J89").Value = 2719.2
$ws.Range("K89").Value = 2787974.5
$ws.Range("L89").Value = 13596
$ws.Range("M89").Value = -2782358.5
$ws.Range("N89").Value = -24828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 385.56
$ws.Range("I94").Value = 341.78262
$ws.Range("J94").Value = 889
$ws.Range("K94").Value = 341.78262
$ws.Range("L94").Value = 889
$ws.Range("M94").Value = 109.21738
$ws.Range("N94").Value = -1791

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4461.875
$ws.Range("I99").Value = 4603.864
$ws.Range("J99").Value = 2900
$ws.Range("K99").Value = 4603.864
$ws.Range("L99").Value = 2900
$ws.Range("M99").Value = -3105.864
$ws.Range("N99").Value = -5896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3102.2424
$ws.Range("I105").Value = 3291.2
$ws.Range("J105").Value = 2811.5386
$ws.Range("K105").Value = 3291.2
$ws.Range("L105").Value = 2811.5386
$ws.Range("M105").Value = -1544.2
$ws.Range("N105").Value = -6305.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2663
$ws.Range("I107").Value = 2663
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2663
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -743
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3000.9092
$ws.Range("J31").Value = 4729.4443
$ws.Range("L31").Value = 4729.4443
$ws.Range("N31").Value = -5319.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3000.9092
$ws.Range("J34").Value = 4729.4443
$ws.Range("L34").Value = 4729.4443
$ws.Range("N34").Value = -5133.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4999.5557
$ws.Range("J86").Value = 4999
$ws.Range("L86").Value = 4999
$ws.Range("N86").Value = -7245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4999.5557
$ws.Range("J89").Value = 4999
$ws.Range("L89").Value = 24995
$ws.Range("N89").Value = -36227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 590.675
$ws.Range("I107").Value = 539.7059
$ws.Range("J107").Value = 879.5
$ws.Range("K107").Value = 539.7059
$ws.Range("L107").Value = 879.5
$ws.Range("M107").Value = 1380.2941
$ws.Range("N107").Value = -4719.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 585.2632
$ws.Range("J12").Value = 977.2727
$ws.Range("L12").Value = 2931.8181
$ws.Range("N12").Value = -3277.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 103.454544
$ws.Range("I33").Value = 143.66667
$ws.Range("J33").Value = 88.375
$ws.Range("K33").Value = 862.0000200000001
$ws.Range("L33").Value = 530.25
$ws.Range("M33").Value = -579.0000200000001
$ws.Range("N33").Value = -1096.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1506.25
$ws.Range("I113").Value = 1430.4
$ws.Range("J113").Value = 1632.6666
$ws.Range("K113").Value = 4291.200000000001
$ws.Range("L113").Value = 4897.9998
$ws.Range("M113").Value = -2121.200000000001
$ws.Range("N113").Value = -9237.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2769.95
$ws.Range("J137").Value = 2741.3333
$ws.Range("L137").Value = 8223.999899999999
$ws.Range("N137").Value = -18423.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 589.8461
$ws.Range("I107").Value = 511.5
$ws.Range("J107").Value = 851
$ws.Range("K107").Value = 511.5
$ws.Range("L107").Value = 851
$ws.Range("M107").Value = 1408.5
$ws.Range("N107").Value = -4691

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 670566.3
$ws.Range("I113").Value = 1999999
$ws.Range("K113").Value = 1999999
$ws.Range("M113").Value = -1997829

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3964.5833
$ws.Range("I132").Value = 3275.111
$ws.Range("K132").Value = 9825.332999999999
$ws.Range("M132").Value = -7295.332999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4150
$ws.Range("I16").Value = 3618
$ws.Range("J16").Value = 10002
$ws.Range("K16").Value = 3618
$ws.Range("L16").Value = 10002
$ws.Range("M16").Value = -3448
$ws.Range("N16").Value = -10342

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4185.273
$ws.Range("I22").Value = 4048.1667
$ws.Range("J22").Value = 4349.8
$ws.Range("K22").Value = 4048.1667
$ws.Range("L22").Value = 4349.8
$ws.Range("M22").Value = -3753.1667
$ws.Range("N22").Value = -4939.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4185.273
$ws.Range("I27").Value = 4048.1667
$ws.Range("J27").Value = 4349.8
$ws.Range("K27").Value = 4048.1667
$ws.Range("L27").Value = 4349.8
$ws.Range("M27").Value = -3941.1667
$ws.Range("N27").Value = -4563.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2079.7856
$ws.Range("I93").Value = 2079.7856
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2079.7856
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -831.7856000000002
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 364.16
$ws.Range("I113").Value = 263.35294
$ws.Range("J113").Value = 578.375
$ws.Range("K113").Value = 790.05882
$ws.Range("L113").Value = 1735.125
$ws.Range("M113").Value = 1379.94118
$ws.Range("N113").Value = -6075.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7657.625
$ws.Range("I122").Value = 7323.143
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 21969.429
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -19519.429
$ws.Range("N122").Value = -34897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10685.625
$ws.Range("I126").Value = 10685.625
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 32056.875
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -29586.875
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19911.674
$ws.Range("I132").Value = 13770
$ws.Range("K132").Value = 41310
$ws.Range("M132").Value = -38780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1485.6154
$ws.Range("I136").Value = 612.6842
$ws.Range("K136").Value = 1838.0526
$ws.Range("M136").Value = 711.9474
